$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Student ID / Staff ID" column (column I) values for rows 2-8.
$ids = @{
    2 = "18-0063"
    3 = "18-0188"
    4 = "18-0023"
    5 = "18-0014"
    6 = "18-0153"
    7 = "18-0253"
    8 = "18-0243"
}

foreach ($row in 2..8) {
    $ws.Range("I$row").Value = $ids[$row]
    # Match the formatting already used by this row's data cells (column B).
    $ws.Range("B$row").Copy()
    $ws.Range("I$row").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# Reflect the new column in the view: select the filled range and zoom out a bit.
$null = $ws.Range("I2:I8").Select()
$excel.ActiveWindow.Zoom = 80
